$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated loading_percent values for rows 2-25 (A column holds line index 0-23)
# Columns changed: B, C, E, F, G, L, N (per-line loading percentages)
$updates = @(
    @{ Row=2; "B"="28.05723335959631"; "C"="11.33315846729864"; "E"="8.67634411646619"; "F"="16.86991607391233"; "G"="3.738699593220593"; "L"="10.8291405045309"; "N"="21.5310185705757" }
    @{ Row=3; "B"="27.57413298666298"; "C"="10.78810392992508"; "E"="8.687263103399868"; "F"="15.89584955866808"; "G"="3.744152897981065"; "L"="10.81529827744079"; "N"="21.55699281753763" }
    @{ Row=4; "B"="27.28292801486936"; "C"="10.44315569470784"; "E"="8.694446387672611"; "F"="15.26997757108491"; "G"="3.747666418550467"; "L"="10.80926756645656"; "N"="21.57477796666122" }
    @{ Row=5; "B"="27.16578143664272"; "C"="10.30021740194889"; "E"="8.697494280514228"; "F"="15.00819731993403"; "G"="3.749139939970195"; "L"="10.80743085387658"; "N"="21.58248510918373" }
    @{ Row=6; "B"="27.14642583914487"; "C"="10.27634646330567"; "E"="8.698007673592258"; "F"="14.96433081551593"; "G"="3.749387143518765"; "L"="10.80716336218253"; "N"="21.58379256156236" }
    @{ Row=7; "B"="27.28134176542476"; "C"="10.441237262165"; "E"="8.694487003852865"; "F"="15.26647399323728"; "G"="3.74768612173236"; "L"="10.80924028224189"; "N"="21.57488005025242" }
    @{ Row=8; "B"="27.88963119174525"; "C"="11.14747827652147"; "E"="8.680009690680675"; "F"="16.53996406344769"; "G"="3.740545732683102"; "L"="10.82385558089159"; "N"="21.53959183847191" }
    @{ Row=9; "B"="29.11844891670251"; "C"="12.44281415377741"; "E"="8.655411423386315"; "F"="19.0027458068253"; "G"="3.727844845111565"; "L"="10.8720791907028"; "N"="21.48507247480618" }
    @{ Row=10; "B"="30.03378230139538"; "C"="13.33101090741333"; "E"="8.639638889705061"; "F"="20.67494806633232"; "G"="3.719294018147316"; "L"="10.91937914800775"; "N"="21.45411428865179" }
    @{ Row=11; "B"="30.45096513170564"; "C"="13.71987963532197"; "E"="8.63296059101193"; "F"="21.3917225636224"; "G"="3.715570736334104"; "L"="10.9434533750155"; "N"="21.44203808138864" }
    @{ Row=12; "B"="30.60889921682782"; "C"="13.86485973852505"; "E"="8.630502950227834"; "F"="21.65686569030329"; "G"="3.714184556653624"; "L"="10.9529347820973"; "N"="21.43775638001252" }
    @{ Row=13; "B"="30.57488969425718"; "C"="13.83373822677926"; "E"="8.631029078829906"; "F"="21.60004134736742"; "G"="3.714482042327583"; "L"="10.95087660492222"; "N"="21.43866552197869" }
    @{ Row=14; "B"="30.4639602294544"; "C"="13.73185332981651"; "E"="8.632756971518775"; "F"="21.4136618050453"; "G"="3.715456219738494"; "L"="10.94422611763069"; "N"="21.44167996960173" }
    @{ Row=15; "B"="30.39600241157134"; "C"="13.66914697887288"; "E"="8.633824635214273"; "F"="21.29868154950795"; "G"="3.716056018033122"; "L"="10.94019995053293"; "N"="21.44356442215556" }
    @{ Row=16; "B"="30.0065207367158"; "C"="13.30528364413457"; "E"="8.640085315105349"; "F"="20.62722412089977"; "G"="3.719540677274127"; "L"="10.9178570775407"; "N"="21.45494411142626" }
    @{ Row=17; "B"="29.76767519012554"; "C"="13.07810620674862"; "E"="8.644053155895511"; "F"="20.20408069617459"; "G"="3.721720913079375"; "L"="10.9048036342494"; "N"="21.46244110479901" }
    @{ Row=18; "B"="29.63038131394408"; "C"="12.9460157344025"; "E"="8.646382113169818"; "F"="19.95656407809808"; "G"="3.722990613903007"; "L"="10.89753652927292"; "N"="21.46694190545882" }
    @{ Row=19; "B"="29.58391526424555"; "C"="12.90105072278338"; "E"="8.647178693585346"; "F"="19.87204792380562"; "G"="3.723423213120024"; "L"="10.89511746071597"; "N"="21.46849813098201" }
    @{ Row=20; "B"="29.79309320065468"; "C"="13.10243774442563"; "E"="8.643625934004115"; "F"="20.2495528364879"; "G"="3.721487201179571"; "L"="10.9061682825414"; "N"="21.46162348478905" }
    @{ Row=21; "B"="30.4965453258165"; "C"="13.76184183104139"; "E"="8.632247513946234"; "F"="21.46857628470567"; "G"="3.715169437352853"; "L"="10.94616964420813"; "N"="21.44078662579298" }
    @{ Row=22; "B"="30.95596665027328"; "C"="14.17950156136102"; "E"="8.625226508787456"; "F"="22.22866616901555"; "G"="3.711178731612653"; "L"="10.97443889467465"; "N"="21.42886770752792" }
    @{ Row=23; "B"="30.71084521402932"; "C"="13.95783240550476"; "E"="8.628935781790259"; "F"="21.82633154475864"; "G"="3.713296056770576"; "L"="10.9591575387378"; "N"="21.43507268049632" }
    @{ Row=24; "B"="29.78160163872663"; "C"="13.09144206860621"; "E"="8.64381893214655"; "F"="20.22900810905294"; "G"="3.721592811685887"; "L"="10.9055505851092"; "N"="21.46199253707782" }
    @{ Row=25; "B"="28.78314905686672"; "C"="12.10299315418301"; "E"="8.661661215493389"; "F"="18.34778573295697"; "G"="3.731142775182205"; "L"="10.85694303537831"; "N"="21.49823457601995" }
)

foreach ($u in $updates) {
    $r = $u.Row
    $ws.Range("B$r").Value = [double]$u["B"]
    $ws.Range("C$r").Value = [double]$u["C"]
    $ws.Range("E$r").Value = [double]$u["E"]
    $ws.Range("F$r").Value = [double]$u["F"]
    $ws.Range("G$r").Value = [double]$u["G"]
    $ws.Range("L$r").Value = [double]$u["L"]
    $ws.Range("N$r").Value = [double]$u["N"]
}
